# Fruta / hortaliza, semanal
#
# Weekly data refresh: the "Fecha" (D) / "Origen" (O) pair for each existing
# observation (rows 34-131, stored as duplicate-row pairs) shifts down by one
# slot to make room for a newly published week at the top of the block
# (rows 34-35), and the pair that falls off the bottom (old rows 130-131)
# is appended as a brand-new pair of rows (132-133) with its original
# Fecha/Origen preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fecha" (serial date) values for each row-pair, top to bottom,
# starting at row 34/35. Index 0 is the freshly published week; indices
# 1..48 are simply the previous pair's old date (i.e. everything shifts
# down by one pair).
$newDates = @(
    44469, 44299, 44320, 44391, 44245, 44334, 44250, 44467, 44434, 44327,
    44383, 44405, 44161, 44336, 44341, 44274, 44442, 44453, 44420, 44237,
    44285, 44344, 44217, 44266, 44350, 44460, 44427, 44280, 44447, 44267,
    44187, 44386, 44308, 44264, 44252, 44166, 44168, 44433, 44371, 44316,
    44209, 44365, 44306, 44203, 44257, 44239, 44376, 44292, 44358
)

# New "Origen" values for each row-pair, same order/shift as above.
$newOrigen = @(
    "Región de Ñuble", "Región de Ñuble", "Región Metropolitana", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región Metropolitana",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble",
    "Región de Ñuble", "Región de Ñuble", "Región de Ñuble", "Región de Ñuble"
)

# The pair of rows pushed off the bottom of the shift (old rows 130/131,
# Fecha = 2021-08-17 / serial 44425) becomes a brand-new pair of rows at
# the end of the sheet, carrying forward the rest of that observation's
# data unchanged. Copy this BEFORE the shift below overwrites row 130/131's
# Fecha/Origen with their new (shifted) values.
$lastOldRow1 = 130
$lastOldRow2 = 131
$newRow1 = 132
$newRow2 = 133

foreach ($pair in @(@($lastOldRow1, $newRow1), @($lastOldRow2, $newRow2))) {
    $srcRow = $pair[0]
    $dstRow = $pair[1]
    for ($col = 1; $col -le 18; $col++) {
        $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
    }
    # Preserve the date number format on column D ("Fecha") for the new row,
    # matching the style used throughout the rest of the column.
    $ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
}

$firstRow = 34
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row1 = $firstRow + (2 * $i)
    $row2 = $row1 + 1

    $ws.Cells.Item($row1, 4).Value = $newDates[$i]
    $ws.Cells.Item($row2, 4).Value = $newDates[$i]

    $ws.Cells.Item($row1, 15).Value = $newOrigen[$i]
    $ws.Cells.Item($row2, 15).Value = $newOrigen[$i]
}
